# Agregué artículo sometido Andrés
# Insert a new supervision entry (Angela Rivero Valderrama y Sebastián Camilo
# Valenzuela) as the first data row of the "supervision" worksheet, pushing
# all the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row right above the current first data row (row 2),
# shifting every row below it down by one.
$ws.Rows(2).Insert()

# Fill in the new entry.
$ws.Range("A2").Value = "Psicología"
$ws.Range("B2").Value = "2023 - 2024"
$ws.Range("C2").Value = "Angela Rivero Valderrama y Sebastián Camilo Valenzuela"
$ws.Range("D2").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E2").Value = "Trabajo de grado: \textit{Preferencias por estímulos sexuales eróticos según género y la orientación sexual: un estudio con eye-tracking}"

# Match the author's final selection/viewport.
$ws.Range("E2").Select()
